# Apply the target edit: in the "Reflection from each member" section, the
# three blank paragraphs that sit directly under each team member's name
# heading ("Damien Tan Lek Khee", "Lim Pau thing", and the final blank
# paragraph that follows "Ooi Ying Jie") get a left indent of 207 twips
# (10.35 points) added to their paragraph properties.

$d = $word.ActiveDocument

function Get-ParagraphAfterHeading($doc, $headingText) {
    $paragraphs = $doc.Paragraphs
    $count = $paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $paragraphs.Item($i)
        $t = $p.Range.Text
        $trimmed = $t.Trim()
        if ($trimmed -eq $headingText -and $p.Style.NameLocal -eq "Heading 2") {
            return $paragraphs.Item($i + 1)
        }
    }
    return $null
}

$targets = @("Damien Tan Lek Khee", "Lim Pau thing", "Ooi Ying Jie")

foreach ($name in $targets) {
    $after = Get-ParagraphAfterHeading $d $name
    if ($after -ne $null) {
        $after.LeftIndent = 10.35
    }
}
